# TC11-Create and Manage API Keys.xlsx — AI Generated sheet used to have two
# columns (A: "button_apiKeyManagement_trNthChild" / "2", B: "input_KeyName" /
# blank). The edit drops column A's old header entirely, keeps only the
# "input_KeyName" header (previously in B1), and removes column B so a
# single column remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B's header ("input_KeyName") becomes the new column A header,
# replacing "button_apiKeyManagement_trNthChild".
$ws.Range("A1").Value = $ws.Range("B1").Value()

# A2 held the placeholder "2"; it becomes blank, same as B2 already was.
# Writing a bare quote-prefix keeps the cell alive as an empty string
# (instead of clearing it away completely), then resetting the style
# drops the implicit "quote prefix" formatting it would otherwise pick up.
$ws.Range("A2").Value = "'"
$ws.Range("A2").Style = "Normal"

# Column A takes on column B's old width (15 character units, i.e. a
# ColumnWidth of 14.17 once Excel's built-in padding is subtracted).
$ws.Range("A1").ColumnWidth = 14.17

# Finally, remove column B so only column A remains.
$ws.Columns.Item(2).Delete()
